$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = '2025-08-25'

$newSheet.Range('A1').Value = 'rank'
$newSheet.Range('B1').Value = 'title'
$newSheet.Range('C1').Value = 'author'
$newSheet.Range('D1').Value = 'latest_episode'

$headerRange = $newSheet.Range('A1:D1')
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$data = New-Object 'object[,]' 50,4
$data[0,0] = 1
$data[0,1] = 'ワンパンマン'
$data[0,2] = '原作/ＯＮＥ 作画/村田雄介'
$data[0,3] = '209撃目'
$data[1,0] = 2
$data[1,1] = '異世界おじさん'
$data[1,2] = '殆ど死んでいる(著者)'
$data[1,3] = '番外編7'
$data[2,0] = 3
$data[2,1] = '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる'
$data[2,2] = '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea'
$data[2,3] = '第70話'
$data[3,0] = 4
$data[3,1] = '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～'
$data[3,2] = 'zunta(作画) はらわたさいぞう(原作)'
$data[3,3] = '第31話：完全なる死角③'
$data[4,0] = 5
$data[4,1] = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$data[4,2] = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$data[4,3] = '第17話-2：「違法奴隷商討伐」'
$data[5,0] = 6
$data[5,1] = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$data[5,2] = '光永康則'
$data[5,3] = '第６８話『施錠停止』③'
$data[6,0] = 7
$data[6,1] = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$data[6,2] = 'マツモトケンゴ'
$data[6,3] = '第６３話　ダンスゲームの戦いが始まった（１）'
$data[7,0] = 8
$data[7,1] = '実は俺、最強でした？'
$data[7,2] = '原作：澄守 彩 漫画：高橋 愛'
$data[7,3] = '第123話　王妃とハルト・前編'
$data[8,0] = 9
$data[8,1] = '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが'
$data[8,2] = '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)'
$data[8,3] = '第4話 後編'
$data[9,0] = 10
$data[9,1] = '元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～'
$data[9,2] = '沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)'
$data[9,3] = '第77話その2'
$data[10,0] = 11
$data[10,1] = 'めっちゃ召喚された件 THE COMIC'
$data[10,2] = '漫画：六甲島カモメ 原作：さいとうさ キャラクター原案：ツグトク'
$data[10,3] = '第47話②'
$data[11,0] = 12
$data[11,1] = '異世界魔王と召喚少女の奴隷魔術'
$data[11,2] = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$data[11,3] = '第127話　戦争を終わらせてみるⅢ（後編）'
$data[12,0] = 13
$data[12,1] = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$data[12,2] = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$data[12,3] = '第５１話　英雄を倒す器用貧乏（４）'
$data[13,0] = 14
$data[13,1] = 'ぽんドロイド！ はまさん'
$data[13,2] = 'はれやまはれぞう(著者)'
$data[13,3] = '第7話'
$data[14,0] = 15
$data[14,1] = 'ダークサモナーとデキている'
$data[14,2] = '車王(著者)'
$data[14,3] = '第74話'
$data[15,0] = 16
$data[15,1] = '不純な彼女達は懺悔しない'
$data[15,2] = 'ポロロッカ(著者)'
$data[15,3] = '第30話'
$data[16,0] = 17
$data[16,1] = '俺の死亡フラグが留まるところを知らない'
$data[16,2] = '漫画：乙須ミツヤ 原作：泉'
$data[16,3] = 'フラグ69 慣れた光景'
$data[17,0] = 18
$data[17,1] = '貧乏騎士に嫁入りしたはずが!? 〜野人令嬢は皇太子妃になっても竜を狩りたい〜'
$data[17,2] = '漫画：夏川そぞろ 原作：宮前葵 キャラクター原案：ののまろ'
$data[17,3] = '第11話④皇族の責務'
$data[18,0] = 19
$data[18,1] = 'まんきつしたい常連さん'
$data[18,2] = 'しんみりん(著者)'
$data[18,3] = '第47話前編'
$data[19,0] = 20
$data[19,1] = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$data[19,2] = '島知宏 音速炒飯 有都あらゆる'
$data[19,3] = '第２３食　巨大ヘビモンスターさん、パクパクですわ！（２）'
$data[20,0] = 21
$data[20,1] = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$data[20,2] = '戸賀 環 坂木持丸 riritto'
$data[20,3] = '第51話①　呪われた家を探索してみた'
$data[21,0] = 22
$data[21,1] = 'リビルドワールド'
$data[21,2] = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$data[21,3] = '第72話③'
$data[22,0] = 23
$data[22,1] = 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ'
$data[22,2] = '板垣恵介 猪原賽 陸井栄史'
$data[22,3] = '第78話　海神(ポセイドン)'
$data[23,0] = 24
$data[23,1] = 'クセ強彼女は床にいざなう'
$data[23,2] = '須河篤志(著者)'
$data[23,3] = '休載イラスト'
$data[24,0] = 25
$data[24,1] = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$data[24,2] = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$data[24,3] = '第33話 独身貴族は見積もりを誤る（1）'
$data[25,0] = 26
$data[25,1] = '異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～'
$data[25,2] = 'モリタ Ｕ４ nima'
$data[25,3] = '第13話（２）　ドーナツの騎士様とお土産スイーツ（２）'
$data[26,0] = 27
$data[26,1] = 'アザミヤコを好きになる'
$data[26,2] = 'ユニティコング(原作) ツノニガウ(作画)'
$data[26,3] = '第9話後編'
$data[27,0] = 28
$data[27,1] = '今日から僕は、彼女の✕✕を解消する。'
$data[27,2] = 'コアヤアコ(著者)'
$data[27,3] = '第2話前半'
$data[28,0] = 29
$data[28,1] = 'ライドンキング'
$data[28,2] = '馬場康誌'
$data[28,3] = '第82話 大統領と星航る龍（前編）'
$data[29,0] = 30
$data[29,1] = '聖者無双'
$data[29,2] = '漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime'
$data[29,3] = '第91話　邂逅（後半）'
$data[30,0] = 31
$data[30,1] = '「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる'
$data[30,2] = '漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨'
$data[30,3] = '再開するよイラスト'
$data[31,0] = 32
$data[31,1] = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$data[31,2] = '村上よしゆき 茨木野 あるてら'
$data[31,3] = '第４１話　勇者、人魚王国を救い、歓迎される。あと、六邪神将が、全員来る（４）'
$data[32,0] = 33
$data[32,1] = '衛宮さんちの今日のごはん'
$data[32,2] = 'TAa(漫画) 只野まこと(料理監修) ＴＹＰＥ－ＭＯＯＮ(原作)'
$data[32,3] = '第75話'
$data[33,0] = 34
$data[33,1] = 'このヒーラー、めんどくさい'
$data[33,2] = '丹念に発酵(著者)'
$data[33,3] = '【コメント募集企画】カーラたちが探検中に転移魔法陣を踏んで飛ばされた先を大募集！【コミックス発売記念】'
$data[34,0] = 35
$data[34,1] = 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～'
$data[34,2] = '板垣恵介 林たかあき'
$data[34,3] = '第52話 洗練されたデザイン'
$data[35,0] = 36
$data[35,1] = '解雇された暗黒兵士(30代)のスローなセカンドライフ'
$data[35,2] = '岡沢六十四 るれくちぇ sage・ジョー'
$data[35,3] = '第72話(前編) ダリエルの照らした道'
$data[36,0] = 37
$data[36,1] = 'ありふれた職業で世界最強'
$data[36,2] = 'RoGa（漫画） 白米 良（原作） たかやKi（キャラクター原案）'
$data[36,3] = '第84話　人間らしさ（前編）'
$data[37,0] = 38
$data[37,1] = 'お気楽領主の楽しい領地防衛 ～生産系魔術で名もなき村を最強の城塞都市に～'
$data[37,2] = '青色まろ（漫画） 赤池宗（原作） 転（原作イラスト）'
$data[37,3] = '第33話　観光案内'
$data[38,0] = 39
$data[38,1] = '異世界黙示録マイノグーラ ～破滅の文明で始める世界征服～'
$data[38,2] = '緑華野菜子(著者) 鹿角フェフ(原作) じゅん(キャラクター原案)'
$data[38,3] = '第31話　帳①'
$data[39,0] = 40
$data[39,1] = '老後に備えて異世界で８万枚の金貨を貯めます'
$data[39,2] = 'FUNA 東西 モトエ恵介'
$data[39,3] = '第121話　会談［その6］'
$data[40,0] = 41
$data[40,1] = 'インフィニット・デンドログラム'
$data[40,2] = '今井神 原作：海道左近 キャラクター原案：タイキ'
$data[40,3] = '第72話'
$data[41,0] = 42
$data[41,1] = 'うちの清楚系委員長がかつて中二病アイドルだったことを俺だけが知っている。'
$data[41,2] = '三上こた こばやし少女 寝子空兄 ゆがー'
$data[41,3] = '第2話　正体'
$data[42,0] = 43
$data[42,1] = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$data[42,2] = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$data[42,3] = '第1話'
$data[43,0] = 44
$data[43,1] = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$data[43,2] = '漫画/すたひろ 原作/Y.A'
$data[43,3] = 'chapter68【36話①】'
$data[44,0] = 45
$data[44,1] = '仮面の黒騎士。正体バレたのでもう学園でも無双する'
$data[44,2] = '楓原こうた(原作) さみ(作画) へいろー(キャラクター原案)'
$data[44,3] = '第1話  恥さらしの正体'
$data[45,0] = 46
$data[45,1] = '魔法歌姫マジカルギンガ'
$data[45,2] = '超銀河レコード'
$data[45,3] = '第14話'
$data[46,0] = 47
$data[46,1] = '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。'
$data[46,2] = '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)'
$data[46,3] = '第24話-2'
$data[47,0] = 48
$data[47,1] = '宇崎ちゃんは遊びたい！'
$data[47,2] = '丈(著者)'
$data[47,3] = '第126話'
$data[48,0] = 49
$data[48,1] = '天獄で悪魔がボクを魅惑する'
$data[48,2] = '銀河味めてお(著者)'
$data[48,3] = '第35話'
$data[49,0] = 50
$data[49,1] = '塔の管理をしてみよう'
$data[49,2] = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$data[49,3] = '第92話前編'

$newSheet.Range('A2:D51').Value = $data

'Sheet 2025-08-25 added with ' + $newSheet.UsedRange.Rows.Count + ' rows'